# "Changes biomarker "TP53" to "p53" label" -- relabel every occurrence of
# the TP53 biomarker in the "translation" sheet's biomarker column.

$wb = $excel.ActiveWorkbook

$wsTranslation  = $wb.Worksheets.Item("translation")
$wsConsolidation = $wb.Worksheets.Item("consolidation")

# Rows 8-17 (col A) hold the "TP53" biomarker rows; relabel them to "p53".
$wsTranslation.Range("A8:A17").Value = "p53"

# Reflect the author's final on-screen view state: the consolidation sheet
# zoomed out and no longer the focused tab, with "translation" active and
# F14 selected there.
$wsConsolidation.Activate()
$excel.ActiveWindow.Zoom = 80
$wsConsolidation.Range("C15").Select() | Out-Null

$wsTranslation.Activate()
$wsTranslation.Range("F14").Select() | Out-Null
